# ECTEST TC25 & TC32 Disabled
# Disable (remove) the sanity rows for TC25_Verify_Footer and
# TC32_Verify_PlaceOrder_ManagedUser from the MasterExecutor sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MasterExecutor")

$lastRow = $ws.UsedRange.Rows.Count
$searchRange = $ws.Range("C1:C$lastRow")

# Delete the TC32 row first (it is below TC25), then TC25, so the
# already-found row number for the first deletion stays valid.
$tc32Cell = $searchRange.Find("TC32_Verify_PlaceOrder_ManagedUser")
if ($tc32Cell -ne $null) {
    $ws.Rows.Item($tc32Cell.Row).Delete()
}

$tc25Cell = $searchRange.Find("TC25_Verify_Footer")
if ($tc25Cell -ne $null) {
    $ws.Rows.Item($tc25Cell.Row).Delete()
}

# Reflect the resulting used range in the selection, matching the
# post-edit column E selection in the authored workbook.
$newLastRow = $ws.UsedRange.Rows.Count
$ws.Range("E1:E$newLastRow").Select()
